$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 97
$ws.Cells.Item($newRow, 1).Value = "22-12-2025"
$ws.Cells.Item($newRow, 2).Value = "The price of gold in India today is ₹13,528 per gram for 24 karat gold, ₹12,400 per gram for 22 karat gold and ₹10,146 per gram for 18 karat gold (also called 999 gold)."

$ws.Cells.Item($newRow, 1).Borders.LineStyle = 1
$ws.Cells.Item($newRow, 2).Borders.LineStyle = 1
$ws.Cells.Item($newRow, 2).WrapText = $true
